$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4985.5293
$ws.Range("I32").Value = 5922.125
$ws.Range("J32").Value = 4153
$ws.Range("K32").Value = 5922.125
$ws.Range("L32").Value = 4153
$ws.Range("M32").Value = -5596.125
$ws.Range("N32").Value = -4805

$ws.Range("H70").Value = 5021.4287
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5021.4287
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 15064.2861
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -15604.2861

$ws.Range("H73").Value = 5021.4287
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5021.4287
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 15064.2861
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -16936.2861

$ws.Range("H75").Value = 57529.4
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 57529.4
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 57529.4
$ws.Range("N75").Value = -59401.4

$ws.Range("H78").Value = 57529.4
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 57529.4
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 172588.2
$ws.Range("N78").Value = -181948.2

$ws.Range("H116").Value = 7062.4717
$ws.Range("I116").Value = 6585.154
$ws.Range("J116").Value = 8392.143
$ws.Range("K116").Value = 6585.154
$ws.Range("L116").Value = 8392.143
$ws.Range("M116").Value = -3143.154
$ws.Range("N116").Value = -15276.143

$ws.Range("H137").Value = 2149.6
$ws.Range("I137").Value = 2140.7646
$ws.Range("J137").Value = 2199.6667
$ws.Range("K137").Value = 6422.293799999999
$ws.Range("L137").Value = 6599.000100000001
$ws.Range("M137").Value = -3872.293799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H45").Value = 1942.72
$ws.Range("I45").Value = 1177.5
$ws.Range("J45").Value = 3910.4285
$ws.Range("K45").Value = 1177.5
$ws.Range("L45").Value = 3910.4285
$ws.Range("M45").Value = -800.5
$ws.Range("N45").Value = -4664.4285

$ws.Range("H51").Value = 46380.332
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 46380.332
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 46380.332
$ws.Range("N51").Value = -47892.332

$ws.Range("H74").Value = 1886.129
$ws.Range("I74").Value = 1730.8422
$ws.Range("J74").Value = 2132
$ws.Range("K74").Value = 1730.8422
$ws.Range("L74").Value = 2132
$ws.Range("M74").Value = -856.8422
$ws.Range("N74").Value = -3880

$ws.Range("H77").Value = 1886.129
$ws.Range("I77").Value = 1730.8422
$ws.Range("J77").Value = 2132
$ws.Range("K77").Value = 8654.210999999999
$ws.Range("L77").Value = 10660
$ws.Range("M77").Value = -4286.210999999999
$ws.Range("N77").Value = -19396

$ws.Range("H80").Value = 33645.5
$ws.Range("I80").Value = 29860.666
$ws.Range("J80").Value = 45000
$ws.Range("K80").Value = 29860.666
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = -28862.666
$ws.Range("N80").Value = -46996

$ws.Range("H83").Value = 33645.5
$ws.Range("I83").Value = 29860.666
$ws.Range("J83").Value = 45000
$ws.Range("K83").Value = 89581.99800000001
$ws.Range("L83").Value = 135000
$ws.Range("M83").Value = -84589.99800000001
$ws.Range("N83").Value = -144984

$ws.Range("H102").Value = 1219.6
$ws.Range("I102").Value = 1115.3125
$ws.Range("J102").Value = 1636.75
$ws.Range("K102").Value = 1115.3125
$ws.Range("L102").Value = 1636.75
$ws.Range("M102").Value = 506.6875
$ws.Range("N102").Value = -4880.75

$ws.Range("H110").Value = 3640.3408
$ws.Range("I110").Value = 3462.4736
$ws.Range("J110").Value = 4766.8335
$ws.Range("K110").Value = 3462.4736
$ws.Range("L110").Value = 4766.8335
$ws.Range("M110").Value = -1417.4736
$ws.Range("N110").Value = -8856.833500000001

$ws.Range("H132").Value = 3221.0667
$ws.Range("I132").Value = 3255.077
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9765.231
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -7235.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 27779082
$ws.Range("I105").Value = 45455668
$ws.Range("J105").Value = 1588.4286
$ws.Range("K105").Value = 45455668
$ws.Range("L105").Value = 1588.4286
$ws.Range("M105").Value = -45453921
$ws.Range("N105").Value = -5082.4286

$ws.Range("H107").Value = 52610.75
$ws.Range("I107").Value = 97818.164
$ws.Range("J107").Value = 7403.3335
$ws.Range("K107").Value = 97818.164
$ws.Range("L107").Value = 7403.3335
$ws.Range("M107").Value = -95898.164
$ws.Range("N107").Value = -11243.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11607.583
$ws.Range("I86").Value = 8447.833000000001
$ws.Range("J86").Value = 14767.333
$ws.Range("K86").Value = 8447.833000000001
$ws.Range("L86").Value = 14767.333
$ws.Range("M86").Value = -7324.833000000001
$ws.Range("N86").Value = -17013.333

$ws.Range("H89").Value = 11607.583
$ws.Range("I89").Value = 8447.833000000001
$ws.Range("J89").Value = 14767.333
$ws.Range("K89").Value = 42239.165
$ws.Range("L89").Value = 73836.66500000001
$ws.Range("M89").Value = -36623.165
$ws.Range("N89").Value = -85068.66500000001

$ws.Range("H105").Value = 1543.3636
$ws.Range("I105").Value = 1693
$ws.Range("J105").Value = 1363.8
$ws.Range("K105").Value = 1693
$ws.Range("L105").Value = 1363.8
$ws.Range("M105").Value = 54
$ws.Range("N105").Value = -4857.8

$ws.Range("H107").Value = 1959.5
$ws.Range("I107").Value = 1680.1
$ws.Range("J107").Value = 3356.5
$ws.Range("K107").Value = 1680.1
$ws.Range("L107").Value = 3356.5
$ws.Range("M107").Value = 239.9000000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 6183.0713
$ws.Range("I64").Value = 2965
$ws.Range("J64").Value = 7060.727
$ws.Range("K64").Value = 8895
$ws.Range("L64").Value = 21182.181
$ws.Range("M64").Value = -8625
$ws.Range("N64").Value = -21722.181

$ws.Range("H67").Value = 6183.0713
$ws.Range("I67").Value = 2965
$ws.Range("J67").Value = 7060.727
$ws.Range("K67").Value = 8895
$ws.Range("L67").Value = 21182.181
$ws.Range("M67").Value = -7959
$ws.Range("N67").Value = -23054.181

$ws.Range("H114").Value = 992
$ws.Range("I114").Value = 988.5
$ws.Range("J114").Value = 999
$ws.Range("K114").Value = 2965.5
$ws.Range("L114").Value = 2997
$ws.Range("M114").Value = 288.5

$ws.Range("H117").Value = 434073.44
$ws.Range("I117").Value = 1503.6
$ws.Range("J117").Value = 630696.0600000001
$ws.Range("K117").Value = 4510.799999999999
$ws.Range("L117").Value = 1892088.18
$ws.Range("M117").Value = -1068.799999999999
$ws.Range("N117").Value = -1898972.18

$ws.Range("H122").Value = 1010.9091
$ws.Range("I122").Value = 599.75
$ws.Range("J122").Value = 1245.8572
$ws.Range("K122").Value = 5397.75
$ws.Range("L122").Value = 11212.7148
$ws.Range("M122").Value = -2947.75
$ws.Range("N122").Value = -16112.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 53514.4
$ws.Range("I36").Value = 102603.8
$ws.Range("J36").Value = 4425
$ws.Range("K36").Value = 102603.8
$ws.Range("L36").Value = 4425
$ws.Range("M36").Value = -102118.8
$ws.Range("N36").Value = -5395

$ws.Range("H80").Value = 4398
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4398
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4398
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6394

$ws.Range("H83").Value = 4398
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4398
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 21990
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -31974

$ws.Range("H122").Value = 61796.75
$ws.Range("I122").Value = 79285.03999999999
$ws.Range("J122").Value = 16327.2
$ws.Range("K122").Value = 237855.12
$ws.Range("L122").Value = 48981.60000000001
$ws.Range("M122").Value = -235405.12

$ws.Range("H132").Value = 4872.4707
$ws.Range("I132").Value = 3253.3333
$ws.Range("J132").Value = 6694
$ws.Range("K132").Value = 9759.999899999999
$ws.Range("L132").Value = 20082
$ws.Range("M132").Value = -7229.999899999999
$ws.Range("N132").Value = -25142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 49999
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 49999
$ws.Range("N42").Value = -51125

$ws.Range("H45").Value = 1697348.5
$ws.Range("I45").Value = 5005999.5
$ws.Range("J45").Value = 43023
$ws.Range("K45").Value = 5005999.5
$ws.Range("L45").Value = 43023
$ws.Range("M45").Value = -5005592.5

$ws.Range("H49").Value = 49999
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 49999
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 49999
$ws.Range("N49").Value = -50293

$ws.Range("H62").Value = 20226
$ws.Range("I62").Value = 20226
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 20226
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -19602

$ws.Range("H65").Value = 20226
$ws.Range("I65").Value = 20226
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 60678
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -57558

$ws.Range("H70").Value = 34747.25
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 34747.25
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 34747.25
$ws.Range("N70").Value = -35287.25

$ws.Range("H73").Value = 34747.25
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 34747.25
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 34747.25
$ws.Range("N73").Value = -36619.25

$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872

$ws.Range("H76").Value = 30288
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30288
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30288
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -30964

$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360

$ws.Range("H79").Value = 30288
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30288
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30288
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -32628

$ws.Range("H136").Value = 4556.207
$ws.Range("I136").Value = 2950.6
$ws.Range("J136").Value = 8124.222
$ws.Range("K136").Value = 8851.799999999999
$ws.Range("L136").Value = 24372.666
$ws.Range("M136").Value = -6301.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 39027.4
$ws.Range("I42").Value = 49999
$ws.Range("J42").Value = 36284.5
$ws.Range("K42").Value = 49999
$ws.Range("L42").Value = 36284.5
$ws.Range("M42").Value = -49621
$ws.Range("N42").Value = -37040.5

$ws.Range("H107").Value = 13285.308
$ws.Range("I107").Value = 2563.625
$ws.Range("J107").Value = 30440
$ws.Range("K107").Value = 7690.875
$ws.Range("L107").Value = 91320
$ws.Range("M107").Value = -5770.875

$ws.Range("H132").Value = 2005.6666
$ws.Range("I132").Value = 2016.5834
$ws.Range("J132").Value = 1874.6666
$ws.Range("K132").Value = 6049.7502
$ws.Range("L132").Value = 5623.9998
$ws.Range("M132").Value = -3519.7502

$ws.Range("H136").Value = 6621.2144
$ws.Range("I136").Value = 8309.888999999999
$ws.Range("J136").Value = 3581.6
$ws.Range("K136").Value = 24929.667
$ws.Range("L136").Value = 10744.8
$ws.Range("M136").Value = -22379.667
